$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "67.947.92"
$ws.Range("E2").Value = "  +1.31%  "

# Row 3
$ws.Range("D3").Value = "3.267.15"
$ws.Range("E3").Value = "  +0.41%  "

# Row 4
$ws.Range("E4").Value = "  -0.01%  "

# Row 5
$ws.Range("D5").Value = "'587.09"
$ws.Range("E5").Value = "  +1.58%  "

# Row 6
$ws.Range("D6").Value = "'186.67"
$ws.Range("E6").Value = "  +4.75%  "

# Row 7
$ws.Range("E7").Value = "  -0.03%  "

# Row 8
$ws.Range("E8").Value = "  -0.41%  "

# Row 9
$ws.Range("E9").Value = "  +3.45%  "

# Row 11
$ws.Range("E11").Value = "  +0.75%  "

# Row 12
$ws.Range("D12").Value = "3.837.26"
$ws.Range("E12").Value = "  +0.55%  "

# Row 13
$ws.Range("E13").Value = "  +0.47%  "

# Row 14
$ws.Range("D14").Value = "'28.53"
$ws.Range("E14").Value = "  +1.27%  "

# Row 15
$ws.Range("D15").Value = "67.960.07"
$ws.Range("E15").Value = "  +1.38%  "

# Row 16
$ws.Range("E16").Value = "  +2.29%  "

# Row 17
$ws.Range("D17").Value = "3.268.76"
$ws.Range("E17").Value = "  +0.39%  "

# Row 18
$ws.Range("D18").Value = "'5.85"
$ws.Range("E18").Value = "  -0.14%  "

# Row 19
$ws.Range("D19").Value = "'13.63"
$ws.Range("E19").Value = "  +1.53%  "

# Row 20
$ws.Range("D20").Value = "'381.68"
$ws.Range("E20").Value = "  +2.23%  "

# Row 21
$ws.Range("E21").Value = "  +0.79%  "

# Row 22
$ws.Range("B22").Value = "Dai"
$ws.Range("C22").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D22").Value = "'0.999"
$ws.Range("E22").Value = "  +0.01%  "

# Row 23
$ws.Range("B23").Value = "Litecoin"
$ws.Range("C23").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D23").Value = "'71.42"
$ws.Range("E23").Value = "  +0.89%  "

# Row 24
$ws.Range("D24").Value = "'0.514"
$ws.Range("E24").Value = "  +0.38%  "

# Row 25
$ws.Range("D25").Value = "'0.0000121"
$ws.Range("E25").Value = "  +1.64%  "

# Row 26
$ws.Range("D26").Value = "'9.83"
$ws.Range("E26").Value = "  -0.84%  "

# Row 27
$ws.Range("E27").Value = "  +4.98%  "

# Row 28
$ws.Range("E28").Value = "  -0.06%  "

# Row 29
$ws.Range("D29").Value = "'5.84"
$ws.Range("E29").Value = "  +3.40%  "

# Row 30
$ws.Range("E30").Value = "  +0.88%  "

# Row 31
$ws.Range("E31").Value = "  +1.06%  "

# Row 32
$ws.Range("D32").Value = "'7.19"
$ws.Range("E32").Value = "  +5.48%  "

# Row 33
$ws.Range("E33").Value = "  -0.01%  "

# Row 34
$ws.Range("E34").Value = "  +0.82%  "

# Row 35
$ws.Range("D35").Value = "'1.53"
$ws.Range("E35").Value = "  +2.09%  "

# Row 36
$ws.Range("D36").Value = "'162.36"
$ws.Range("E36").Value = "  -2.53%  "

# Row 37
$ws.Range("E37").Value = "  -1.02%  "

# Row 38
$ws.Range("D38").Value = "'0.837"
$ws.Range("E38").Value = "  -2.20%  "

# Row 39
$ws.Range("D39").Value = "'6.78"
$ws.Range("E39").Value = "  +4.42%  "

# Row 40
$ws.Range("D40").Value = "'26.53"
$ws.Range("E40").Value = "  -2.07%  "

# Row 41
$ws.Range("E41").Value = "  +4.49%  "

# Row 42
$ws.Range("D42").Value = "'2.62"
$ws.Range("E42").Value = "  +0.81%  "

# Row 43
$ws.Range("B43").Value = "Hedera"
$ws.Range("C43").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D43").Value = "'0.0689"
$ws.Range("E43").Value = "  +2.08%  "

# Row 44
$ws.Range("B44").Value = "OKB"
$ws.Range("C44").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D44").Value = "'41.20"
$ws.Range("E44").Value = "  +1.77%  "

# Row 45
$ws.Range("D45").Value = "'25.41"
$ws.Range("E45").Value = "  -0.27%  "

# Row 46
$ws.Range("D46").Value = "2.641.37"
$ws.Range("E46").Value = "  -4.32%  "

# Row 47
$ws.Range("D47").Value = "'343.54"
$ws.Range("E47").Value = "  -2.91%  "

# Row 48
$ws.Range("D48").Value = "'0.0285"
$ws.Range("E48").Value = "  +1.59%  "

# Row 49
$ws.Range("D49").Value = "'32.06"
$ws.Range("E49").Value = "  +3.78%  "

# Row 50
$ws.Range("E50").Value = "  +0.97%  "

# Row 51
$ws.Range("E51").Value = "  -0.60%  "
